$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append after the existing data (row 55)
# Date values are Excel serial date numbers; copy formatting/style from the
# last existing data row (A55) so the new date cells inherit the same
# number format (m/d/yyyy) used throughout column A.
$ws.Range("A55").Copy($ws.Range("A56:A59"))

$ws.Range("A56").Value = 46031
$ws.Range("B56").Value = 15

$ws.Range("A57").Value = 46036
$ws.Range("B57").Value = 5

$ws.Range("A58").Value = 46034
$ws.Range("B58").Value = 1

$ws.Range("A59").Value = 46035
$ws.Range("B59").Value = 3

# Update selection to match the author's saved view state
$ws.Range("A56:B56").Select()
